$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing value in B94 ---
$ws.Range("B94").Value2 = 129.867934141357

# --- Add new row 95 data ---

# Column A: date serial; copy A94's cell (style/number format) down to A95 first
# so the new date cell picks up the same built-in date style, then set its value.
$ws.Range("A94").Copy($ws.Range("A95"))
$ws.Range("A95").Value2 = 45566

# Columns B & C: plain numeric values.
$ws.Range("B95").Value2 = 105.352990696857
$ws.Range("C95").Value2 = 121.48502296079

# Columns D-G: text values that look numeric (must stay text, not become numbers).
# Build each as a formula returning the literal string, then paste back as a value
# only; this converts the cell to a genuine text cell without leaving behind any
# extra/leftover number-format style on the cell.
$ws.Range("D95").Formula = '="110.4"'
$ws.Range("D95").Copy()
$ws.Range("D95").PasteSpecial(-4163)

$ws.Range("E95").Formula = '="112.1"'
$ws.Range("E95").Copy()
$ws.Range("E95").PasteSpecial(-4163)

$ws.Range("F95").Formula = '=" 88.5"'
$ws.Range("F95").Copy()
$ws.Range("F95").PasteSpecial(-4163)

$ws.Range("G95").Formula = '="169.0"'
$ws.Range("G95").Copy()
$ws.Range("G95").PasteSpecial(-4163)
